# Auto-generated edit script: updates crypto price (D) and volume-change (E) columns
# per the commit "Updated cryptos list on Fri Jun 16 23:07:33 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving it as TEXT (matching the
# source file, which stores every D/E cell as an inline string) even when the
# text looks like a plain number (Excel's COM layer would otherwise silently
# coerce "0.9996" etc. into a numeric cell). Forcing the text number format,
# assigning the value, then clearing the formatting again keeps the final cell
# style identical to the original (no explicit style index).
function Set-TextValue($cell, $val) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "26.355.67"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "1.718.75"
$ws.Range("E3").Value = "  +3.27%  "
Set-TextValue "D5" "239.14"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("E6").Value = "  -0.01%  "
Set-TextValue "D7" "0.4734"
$ws.Range("E7").Value = "  -1.27%  "
Set-TextValue "D8" "0.2634"
$ws.Range("E8").Value = "  +0.58%  "
Set-TextValue "D9" "0.06210"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "1.716.28"
$ws.Range("E10").Value = "  +3.25%  "
Set-TextValue "D11" "0.07074"
$ws.Range("E11").Value = "  -0.15%  "
Set-TextValue "D12" "15.32"
$ws.Range("E12").Value = "  +3.83%  "
Set-TextValue "D13" "0.5911"
$ws.Range("E13").Value = "  +0.00%  "
Set-TextValue "D14" "4.417"
$ws.Range("E14").Value = "  +0.83%  "
Set-TextValue "D15" "76.20"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "26.346.36"
$ws.Range("E18").Value = "  +2.97%  "
Set-TextValue "D19" "0.000006811"
$ws.Range("E19").Value = "  +0.85%  "
Set-TextValue "D20" "11.54"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D21").Value = "1.936.06"
$ws.Range("E21").Value = "  +3.17%  "
Set-TextValue "D22" "4.543"
$ws.Range("E22").Value = "  +2.54%  "
Set-TextValue "D23" "8.764"
$ws.Range("E23").Value = "  +1.29%  "
Set-TextValue "D24" "5.322"
$ws.Range("E24").Value = "  +0.20%  "
Set-TextValue "D25" "134.75"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  -0.02%  "
Set-TextValue "D28" "108.04"
$ws.Range("E28").Value = "  +3.21%  "
Set-TextValue "D29" "1.755"
$ws.Range("E29").Value = "  +3.96%  "
Set-TextValue "D30" "4.003"
$ws.Range("E30").Value = "  +1.38%  "
Set-TextValue "D31" "3.690"
$ws.Range("E31").Value = "  +0.86%  "
Set-TextValue "D32" "0.07745"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +2.53%  "
Set-TextValue "D34" "2.611"
$ws.Range("E34").Value = "  -0.34%  "
Set-TextValue "D35" "0.9770"
$ws.Range("E35").Value = "  +2.77%  "
Set-TextValue "D36" "0.6193"
$ws.Range("E36").Value = "  +1.19%  "
Set-TextValue "D37" "0.9376"
$ws.Range("E37").Value = "  +9.58%  "
Set-TextValue "D38" "113.95"
$ws.Range("E38").Value = "  +16.30%  "
Set-TextValue "D39" "2.413"
$ws.Range("E39").Value = "  -7.53%  "
Set-TextValue "D40" "1.922"
$ws.Range("E40").Value = "  +2.28%  "
Set-TextValue "D41" "0.9996"
$ws.Range("E41").Value = "  -0.05%  "
Set-TextValue "D42" "0.01477"
$ws.Range("E42").Value = "  -1.69%  "
Set-TextValue "D43" "5.333"
$ws.Range("E43").Value = "  +13.37%  "
Set-TextValue "D44" "0.3817"
$ws.Range("E44").Value = "  +1.39%  "
Set-TextValue "D45" "0.1170"
$ws.Range("E45").Value = "  +4.46%  "
Set-TextValue "D46" "6.287"
$ws.Range("E46").Value = "  +1.18%  "
Set-TextValue "D47" "0.05283"
$ws.Range("E47").Value = "  +0.36%  "
Set-TextValue "D48" "30.35"
$ws.Range("E48").Value = "  +2.94%  "
Set-TextValue "D49" "7.722"
$ws.Range("E49").Value = "  +5.51%  "
Set-TextValue "D50" "1.216"
$ws.Range("E50").Value = "  +1.54%  "
Set-TextValue "D51" "0.3368"
$ws.Range("E51").Value = "  +0.96%  "
